# Re-sort the calibration data (rows 2-8) in ascending order by column A (time).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @()
for ($r = 2; $r -le 8; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $d = $ws.Cells.Item($r, 4).Value()
    $rows += ,@($a, $b, $c, $d)
}

$sorted = $rows | Sort-Object { $_[0] }

for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $sorted[$i][0]
    $ws.Cells.Item($r, 2).Value = $sorted[$i][1]
    $ws.Cells.Item($r, 3).Value = $sorted[$i][2]
    $ws.Cells.Item($r, 4).Value = $sorted[$i][3]
}
